# Updated cryptos list on Sun Oct 29 17:36:53 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for the crypto
# table with the latest scrape, and fixes a rank-13/14 ordering swap
# (WrappedEther now outranks Chainlink) including its Coin name (B) and
# Link (C) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Force 'Price' (column D) cells that are being updated to Text format ---
# The source workbook stores these as inline strings (t="inlineStr"), since the
# upstream price feed uses locale-style thousands separators (e.g. 34.486.86) that
# are not valid Excel numbers. Several of the new values (e.g. 0.999, 33.03) WOULD
# be auto-coerced to numeric by a plain .Value assignment, so pre-set NumberFormat
# to Text ("@") on each touched Price cell to keep them as text, matching the feed.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# --- Apply the updated cell values row by row ---
$ws.Range("D2").Value = '34.486.86'
$ws.Range("E2").Value = '  +0.94%  '

$ws.Range("D3").Value = '1.794.39'
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").Value = '226.82'
$ws.Range("E5").Value = '  +0.11%  '

$ws.Range("E6").Value = '  +2.00%  '

$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("D8").Value = '33.03'
$ws.Range("E8").Value = '  +3.80%  '

$ws.Range("D9").Value = '0.298'
$ws.Range("E9").Value = '  +1.97%  '

$ws.Range("D10").Value = '0.0695'
$ws.Range("E10").Value = '  +0.87%  '

$ws.Range("E11").Value = '  +0.44%  '

$ws.Range("D12").Value = '2.050.91'
$ws.Range("E12").Value = '  +0.24%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.862.47'
$ws.Range("E13").Value = '  +4.14%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '11.15'
$ws.Range("E14").Value = '  +1.07%  '

$ws.Range("D15").Value = '0.638'
$ws.Range("E15").Value = '  +2.18%  '

$ws.Range("D16").Value = '34.396.49'
$ws.Range("E16").Value = '  +0.82%  '

$ws.Range("D17").Value = '4.28'
$ws.Range("E17").Value = '  +2.41%  '

$ws.Range("D18").Value = '68.88'
$ws.Range("E18").Value = '  +0.90%  '

$ws.Range("D19").Value = '248.86'
$ws.Range("E19").Value = '  +0.71%  '

$ws.Range("D20").Value = '0.0₃0801'
$ws.Range("E20").Value = '  +3.08%  '

$ws.Range("D21").Value = '11.38'
$ws.Range("E21").Value = '  +4.20%  '

$ws.Range("D23").Value = '4.17'
$ws.Range("E23").Value = '  +1.66%  '

$ws.Range("E24").Value = '  +1.36%  '

$ws.Range("D25").Value = '165.06'
$ws.Range("E25").Value = '  +2.43%  '

$ws.Range("D26").Value = '7.27'
$ws.Range("E26").Value = '  +1.02%  '

$ws.Range("D27").Value = '16.57'
$ws.Range("E27").Value = '  +1.35%  '

$ws.Range("E28").Value = '  +2.83%  '

$ws.Range("E29").Value = '  -0.16%  '

$ws.Range("D30").Value = '3.81'
$ws.Range("E30").Value = '  +3.22%  '

$ws.Range("D31").Value = '3.92'
$ws.Range("E31").Value = '  +7.88%  '

$ws.Range("E32").Value = '  -0.16%  '

$ws.Range("E33").Value = '  +0.31%  '

$ws.Range("D34").Value = '1.83'
$ws.Range("E34").Value = '  +1.75%  '

$ws.Range("D35").Value = '1.419.83'
$ws.Range("E35").Value = '  -1.81%  '

$ws.Range("D36").Value = '2.60'
$ws.Range("E36").Value = '  +5.72%  '

$ws.Range("D37").Value = '0.675'
$ws.Range("E37").Value = '  +3.23%  '

$ws.Range("E38").Value = '  +0.46%  '

$ws.Range("E39").Value = '  +1.70%  '

$ws.Range("D40").Value = '85.09'
$ws.Range("E40").Value = '  +5.54%  '

$ws.Range("D41").Value = '2.39'
$ws.Range("E41").Value = '  +0.80%  '

$ws.Range("D42").Value = '0.939'
$ws.Range("E42").Value = '  +1.63%  '

$ws.Range("D43").Value = '2.74'
$ws.Range("E43").Value = '  +1.99%  '

$ws.Range("D44").Value = '13.50'
$ws.Range("E44").Value = '  -0.31%  '

$ws.Range("D45").Value = '0.0522'
$ws.Range("E45").Value = '  +2.69%  '

$ws.Range("D46").Value = '6.10'
$ws.Range("E46").Value = '  +0.45%  '

$ws.Range("E47").Value = '  -0.12%  '

$ws.Range("D48").Value = '1.950.70'
$ws.Range("E48").Value = '  +0.15%  '

$ws.Range("D49").Value = '105.75'
$ws.Range("E49").Value = '  -0.32%  '

$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("E51").Value = '  -5.56%  '

